# Update "想去人数" (want-to-go count) values in column F for the sheets
# that hold the full exhibition listing: "展览" and "全部类型".
# Both sheets carry the same rows/data, so the same cell updates apply to each.

$wb = $excel.ActiveWorkbook

$updates = @{
    4  = 47
    6  = 3002
    8  = 2058
    11 = 903
    14 = 227
    16 = 94
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
